$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 1.009582666666667
$ws.Range("H2").Value = 3.028748
$ws.Range("I2").Value = 0.2254436189979109
$ws.Range("J2").Value = 0.2254436189979109
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 26.532132
$ws.Range("N2").Value = 79.596396
$ws.Range("O2").Value = 0.3960736634233649
$ws.Range("P2").Value = 0.3960736634233648
$ws.Range("Q2").Value = 26.786380576912
$ws.Range("R2").Value = 241.077425192208
$ws.Range("S2").Value = 0.08929228007192386
$ws.Range("T2").Value = 0.08929228007192386
$ws.Range("G3").Value = 1.009582666666667
$ws.Range("H3").Value = 3.028748
$ws.Range("I3").Value = 0.2254436189979109
$ws.Range("J3").Value = 0.2254436189979109
$ws.Range("O3").Value = 0.2505213219764053
$ws.Range("P3").Value = 0.2505213219764053
$ws.Range("Q3").Value = 16.94270559443422
$ws.Range("R3").Value = 152.484350349908
$ws.Range("S3").Value = 0.05647843346250168
$ws.Range("T3").Value = 0.05647843346250168
$ws.Range("G4").Value = 1.009582666666667
$ws.Range("H4").Value = 3.028748
$ws.Range("I4").Value = 0.2254436189979109
$ws.Range("J4").Value = 0.2254436189979109
$ws.Range("M4").Value = 23.67385
$ws.Range("N4").Value = 71.02154999999999
$ws.Range("O4").Value = 0.3534050146002298
$ws.Range("P4").Value = 0.3534050146002298
$ws.Range("Q4").Value = 23.90070861326666
$ws.Range("R4").Value = 215.1063775194
$ws.Range("S4").Value = 0.07967290546348535
$ws.Range("T4").Value = 0.07967290546348535
$ws.Range("I5").Value = 0.4390905462561113
$ws.Range("J5").Value = 0.4390905462561113
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 26.532132
$ws.Range("N5").Value = 79.596396
$ws.Range("O5").Value = 0.3960736634233649
$ws.Range("P5").Value = 0.3960736634233648
$ws.Range("Q5").Value = 52.171121684528
$ws.Range("R5").Value = 469.540095160752
$ws.Range("S5").Value = 0.1739122012302244
$ws.Range("T5").Value = 0.1739122012302244
$ws.Range("I6").Value = 0.4390905462561113
$ws.Range("J6").Value = 0.4390905462561113
$ws.Range("O6").Value = 0.2505213219764053
$ws.Range("P6").Value = 0.2505213219764053
$ws.Range("S6").Value = 0.1100015441154229
$ws.Range("T6").Value = 0.1100015441154229
$ws.Range("I7").Value = 0.4390905462561113
$ws.Range("J7").Value = 0.4390905462561113
$ws.Range("M7").Value = 23.67385
$ws.Range("N7").Value = 71.02154999999999
$ws.Range("O7").Value = 0.3534050146002298
$ws.Range("P7").Value = 0.3534050146002298
$ws.Range("Q7").Value = 46.55077507873333
$ws.Range("R7").Value = 418.9569757085999
$ws.Range("S7").Value = 0.1551768009104639
$ws.Range("T7").Value = 0.1551768009104639
$ws.Range("G8").Value = 1.502284666666667
$ws.Range("H8").Value = 4.506854000000001
$ws.Range("I8").Value = 0.3354658347459779
$ws.Range("J8").Value = 0.3354658347459779
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 26.532132
$ws.Range("N8").Value = 79.596396
$ws.Range("O8").Value = 0.3960736634233649
$ws.Range("P8").Value = 0.3960736634233648
$ws.Range("Q8").Value = 39.85881507757601
$ws.Range("R8").Value = 358.7293356981841
$ws.Range("S8").Value = 0.1328691821212166
$ws.Range("T8").Value = 0.1328691821212165
$ws.Range("G9").Value = 1.502284666666667
$ws.Range("H9").Value = 4.506854000000001
$ws.Range("I9").Value = 0.3354658347459779
$ws.Range("J9").Value = 0.3354658347459779
$ws.Range("O9").Value = 0.2505213219764053
$ws.Range("P9").Value = 0.2505213219764053
$ws.Range("Q9").Value = 25.21117652544823
$ws.Range("R9").Value = 226.900588729034
$ws.Range("S9").Value = 0.0840413443984807
$ws.Range("T9").Value = 0.0840413443984807
$ws.Range("G10").Value = 1.502284666666667
$ws.Range("H10").Value = 4.506854000000001
$ws.Range("I10").Value = 0.3354658347459779
$ws.Range("J10").Value = 0.3354658347459779
$ws.Range("M10").Value = 23.67385
$ws.Range("N10").Value = 71.02154999999999
$ws.Range("O10").Value = 0.3534050146002298
$ws.Range("P10").Value = 0.3534050146002298
$ws.Range("Q10").Value = 35.56486185596667
$ws.Range("R10").Value = 320.0837567037
$ws.Range("S10").Value = 0.1185553082262806
$ws.Range("T10").Value = 0.1185553082262806
